# Adds "longest" and "shortest" example-length stats to the basic stats table.
# Inserts 6 new columns (train_longest, dev_longest, test_longest,
# train_shortest, dev_shortest, test_shortest) right after test_avg_tokens
# (i.e. before the old train_hapaxes column), shifting the hapax/unknown
# columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank columns at H:M (everything from the old column H onward
# shifts right by 6 columns, e.g. old H -> new N, old S -> new Y).
$ws.Range("H1:M1").EntireColumn.Insert()

# New header row labels for the inserted columns. The inserted columns
# already inherit the header row's (bold/centered/bordered) style from the
# EntireColumn.Insert() above, so there's no need to re-apply style here
# (doing so via Range.Style would reset it to the default style instead).
$ws.Range("H1").Value2 = "train_longest"
$ws.Range("I1").Value2 = "dev_longest"
$ws.Range("J1").Value2 = "test_longest"
$ws.Range("K1").Value2 = "train_shortest"
$ws.Range("L1").Value2 = "dev_shortest"
$ws.Range("M1").Value2 = "test_shortest"

# New per-language data: train_longest, dev_longest, test_longest,
# train_shortest, dev_shortest, test_shortest.
$data = @{
    2  = @(3819, 3367, 3089, 6, 8, 6)
    3  = @(166, 104, 169, 4, 4, 4)
    4  = @(404, 325, 411, 3, 3, 3)
    5  = @(140, 85, 91, 3, 4, 5)
    6  = @(2254, 495, 442, 6, 15, 6)
    7  = @(2878, 2791, 847, 3, 4, 4)
    8  = @(478, 476, 826, 3, 5, 4)
    9  = @(1103, 999, 1336, 3, 3, 3)
    10 = @(250, 250, 250, 3, 6, 3)
    11 = @(176, 140, 146, 3, 3, 3)
    12 = @(3276, 671, 1254, 3, 4, 3)
    13 = @(146, 74, 79, 5, 5, 5)
    14 = @(142, 114, 112, 3, 3, 3)
    15 = @(866, 789, 767, 9, 11, 10)
    16 = @(635, 67, 122, 3, 6, 7)
    17 = @(406, 148, 957, 3, 3, 4)
    18 = @(2224, 1549, 1169, 3, 4, 3)
    19 = @(177, 114, 107, 6, 8, 6)
    20 = @(678, 257, 846, 6, 11, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 8 + $i).Value2 = $vals[$i]
    }
}
